$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Samurai (5)")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Samurai (6)"

$grid = New-Object 'object[,]' 21,21
$grid[0,0] = ""
$grid[0,1] = ""
$grid[0,2] = 3
$grid[0,3] = ""
$grid[0,4] = 9
$grid[0,5] = ""
$grid[0,6] = 5
$grid[0,7] = ""
$grid[0,8] = ""
$grid[0,9] = ""
$grid[0,10] = ""
$grid[0,11] = ""
$grid[0,12] = 1
$grid[0,13] = ""
$grid[0,14] = ""
$grid[0,15] = 5
$grid[0,16] = ""
$grid[0,17] = ""
$grid[0,18] = ""
$grid[0,19] = 7
$grid[0,20] = 3
$grid[1,0] = ""
$grid[1,1] = 5
$grid[1,2] = ""
$grid[1,3] = ""
$grid[1,4] = ""
$grid[1,5] = 4
$grid[1,6] = ""
$grid[1,7] = 8
$grid[1,8] = ""
$grid[1,9] = ""
$grid[1,10] = ""
$grid[1,11] = ""
$grid[1,12] = ""
$grid[1,13] = ""
$grid[1,14] = ""
$grid[1,15] = 4
$grid[1,16] = ""
$grid[1,17] = 3
$grid[1,18] = ""
$grid[1,19] = ""
$grid[1,20] = ""
$grid[2,0] = ""
$grid[2,1] = 7
$grid[2,2] = ""
$grid[2,3] = ""
$grid[2,4] = ""
$grid[2,5] = ""
$grid[2,6] = ""
$grid[2,7] = ""
$grid[2,8] = 9
$grid[2,9] = ""
$grid[2,10] = ""
$grid[2,11] = ""
$grid[2,12] = ""
$grid[2,13] = 2
$grid[2,14] = ""
$grid[2,15] = ""
$grid[2,16] = ""
$grid[2,17] = ""
$grid[2,18] = ""
$grid[2,19] = 5
$grid[2,20] = 1
$grid[3,0] = ""
$grid[3,1] = 6
$grid[3,2] = ""
$grid[3,3] = ""
$grid[3,4] = 2
$grid[3,5] = 5
$grid[3,6] = ""
$grid[3,7] = 4
$grid[3,8] = ""
$grid[3,9] = ""
$grid[3,10] = ""
$grid[3,11] = ""
$grid[3,12] = ""
$grid[3,13] = ""
$grid[3,14] = 5
$grid[3,15] = 7
$grid[3,16] = ""
$grid[3,17] = 4
$grid[3,18] = 9
$grid[3,19] = ""
$grid[3,20] = ""
$grid[4,0] = ""
$grid[4,1] = ""
$grid[4,2] = 6
$grid[4,3] = ""
$grid[4,4] = ""
$grid[4,5] = ""
$grid[4,6] = 8
$grid[4,7] = ""
$grid[4,8] = ""
$grid[4,9] = ""
$grid[4,10] = ""
$grid[4,11] = ""
$grid[4,12] = ""
$grid[4,13] = ""
$grid[4,14] = ""
$grid[4,15] = ""
$grid[4,16] = ""
$grid[4,17] = ""
$grid[4,18] = ""
$grid[4,19] = ""
$grid[4,20] = ""
$grid[5,0] = ""
$grid[5,1] = 3
$grid[5,2] = ""
$grid[5,3] = 1
$grid[5,4] = 4
$grid[5,5] = ""
$grid[5,6] = ""
$grid[5,7] = 2
$grid[5,8] = ""
$grid[5,9] = ""
$grid[5,10] = ""
$grid[5,11] = ""
$grid[5,12] = ""
$grid[5,13] = ""
$grid[5,14] = 9
$grid[5,15] = 2
$grid[5,16] = ""
$grid[5,17] = 1
$grid[5,18] = 6
$grid[5,19] = ""
$grid[5,20] = ""
$grid[6,0] = 3
$grid[6,1] = ""
$grid[6,2] = ""
$grid[6,3] = ""
$grid[6,4] = ""
$grid[6,5] = ""
$grid[6,6] = ""
$grid[6,7] = ""
$grid[6,8] = ""
$grid[6,9] = ""
$grid[6,10] = 6
$grid[6,11] = ""
$grid[6,12] = ""
$grid[6,13] = 3
$grid[6,14] = ""
$grid[6,15] = ""
$grid[6,16] = ""
$grid[6,17] = ""
$grid[6,18] = ""
$grid[6,19] = 2
$grid[6,20] = ""
$grid[7,0] = ""
$grid[7,1] = 2
$grid[7,2] = ""
$grid[7,3] = 7
$grid[7,4] = ""
$grid[7,5] = ""
$grid[7,6] = 6
$grid[7,7] = ""
$grid[7,8] = ""
$grid[7,9] = 4
$grid[7,10] = ""
$grid[7,11] = 9
$grid[7,12] = ""
$grid[7,13] = ""
$grid[7,14] = ""
$grid[7,15] = 6
$grid[7,16] = ""
$grid[7,17] = 5
$grid[7,18] = ""
$grid[7,19] = ""
$grid[7,20] = ""
$grid[8,0] = ""
$grid[8,1] = ""
$grid[8,2] = 7
$grid[8,3] = ""
$grid[8,4] = 1
$grid[8,5] = ""
$grid[8,6] = 3
$grid[8,7] = 9
$grid[8,8] = ""
$grid[8,9] = ""
$grid[8,10] = ""
$grid[8,11] = ""
$grid[8,12] = ""
$grid[8,13] = 6
$grid[8,14] = 8
$grid[8,15] = ""
$grid[8,16] = ""
$grid[8,17] = 9
$grid[8,18] = ""
$grid[8,19] = ""
$grid[8,20] = 4
$grid[9,0] = ""
$grid[9,1] = ""
$grid[9,2] = ""
$grid[9,3] = ""
$grid[9,4] = ""
$grid[9,5] = ""
$grid[9,6] = ""
$grid[9,7] = ""
$grid[9,8] = ""
$grid[9,9] = 8
$grid[9,10] = ""
$grid[9,11] = 5
$grid[9,12] = ""
$grid[9,13] = ""
$grid[9,14] = 9
$grid[9,15] = ""
$grid[9,16] = ""
$grid[9,17] = ""
$grid[9,18] = ""
$grid[9,19] = ""
$grid[9,20] = ""
$grid[10,0] = ""
$grid[10,1] = ""
$grid[10,2] = ""
$grid[10,3] = ""
$grid[10,4] = ""
$grid[10,5] = ""
$grid[10,6] = ""
$grid[10,7] = 3
$grid[10,8] = ""
$grid[10,9] = ""
$grid[10,10] = 9
$grid[10,11] = ""
$grid[10,12] = ""
$grid[10,13] = 1
$grid[10,14] = ""
$grid[10,15] = ""
$grid[10,16] = ""
$grid[10,17] = ""
$grid[10,18] = ""
$grid[10,19] = ""
$grid[10,20] = ""
$grid[11,0] = ""
$grid[11,1] = ""
$grid[11,2] = ""
$grid[11,3] = ""
$grid[11,4] = ""
$grid[11,5] = ""
$grid[11,6] = 2
$grid[11,7] = ""
$grid[11,8] = ""
$grid[11,9] = 3
$grid[11,10] = ""
$grid[11,11] = 6
$grid[11,12] = ""
$grid[11,13] = ""
$grid[11,14] = ""
$grid[11,15] = ""
$grid[11,16] = ""
$grid[11,17] = ""
$grid[11,18] = ""
$grid[11,19] = ""
$grid[11,20] = ""
$grid[12,0] = ""
$grid[12,1] = ""
$grid[12,2] = ""
$grid[12,3] = ""
$grid[12,4] = ""
$grid[12,5] = ""
$grid[12,6] = 9
$grid[12,7] = 8
$grid[12,8] = ""
$grid[12,9] = ""
$grid[12,10] = ""
$grid[12,11] = ""
$grid[12,12] = ""
$grid[12,13] = 5
$grid[12,14] = 7
$grid[12,15] = ""
$grid[12,16] = ""
$grid[12,17] = ""
$grid[12,18] = 9
$grid[12,19] = 2
$grid[12,20] = ""
$grid[13,0] = 1
$grid[13,1] = ""
$grid[13,2] = ""
$grid[13,3] = ""
$grid[13,4] = 6
$grid[13,5] = ""
$grid[13,6] = ""
$grid[13,7] = ""
$grid[13,8] = ""
$grid[13,9] = 9
$grid[13,10] = ""
$grid[13,11] = 1
$grid[13,12] = ""
$grid[13,13] = ""
$grid[13,14] = 6
$grid[13,15] = 5
$grid[13,16] = ""
$grid[13,17] = 9
$grid[13,18] = 3
$grid[13,19] = ""
$grid[13,20] = ""
$grid[14,0] = ""
$grid[14,1] = 6
$grid[14,2] = 5
$grid[14,3] = ""
$grid[14,4] = 9
$grid[14,5] = ""
$grid[14,6] = ""
$grid[14,7] = 2
$grid[14,8] = ""
$grid[14,9] = ""
$grid[14,10] = 5
$grid[14,11] = ""
$grid[14,12] = ""
$grid[14,13] = ""
$grid[14,14] = ""
$grid[14,15] = 2
$grid[14,16] = ""
$grid[14,17] = ""
$grid[14,18] = ""
$grid[14,19] = ""
$grid[14,20] = ""
$grid[15,0] = ""
$grid[15,1] = 4
$grid[15,2] = ""
$grid[15,3] = 6
$grid[15,4] = ""
$grid[15,5] = 9
$grid[15,6] = ""
$grid[15,7] = 1
$grid[15,8] = ""
$grid[15,9] = ""
$grid[15,10] = ""
$grid[15,11] = ""
$grid[15,12] = 3
$grid[15,13] = ""
$grid[15,14] = ""
$grid[15,15] = ""
$grid[15,16] = ""
$grid[15,17] = ""
$grid[15,18] = ""
$grid[15,19] = ""
$grid[15,20] = 5
$grid[16,0] = ""
$grid[16,1] = 8
$grid[16,2] = ""
$grid[16,3] = ""
$grid[16,4] = ""
$grid[16,5] = ""
$grid[16,6] = ""
$grid[16,7] = 4
$grid[16,8] = ""
$grid[16,9] = ""
$grid[16,10] = ""
$grid[16,11] = ""
$grid[16,12] = ""
$grid[16,13] = ""
$grid[16,14] = 9
$grid[16,15] = 4
$grid[16,16] = ""
$grid[16,17] = 8
$grid[16,18] = 2
$grid[16,19] = ""
$grid[16,20] = ""
$grid[17,0] = ""
$grid[17,1] = 5
$grid[17,2] = ""
$grid[17,3] = 1
$grid[17,4] = ""
$grid[17,5] = 2
$grid[17,6] = ""
$grid[17,7] = 3
$grid[17,8] = ""
$grid[17,9] = ""
$grid[17,10] = ""
$grid[17,11] = ""
$grid[17,12] = 2
$grid[17,13] = ""
$grid[17,14] = ""
$grid[17,15] = ""
$grid[17,16] = ""
$grid[17,17] = ""
$grid[17,18] = ""
$grid[17,19] = ""
$grid[17,20] = 7
$grid[18,0] = ""
$grid[18,1] = 3
$grid[18,2] = ""
$grid[18,3] = ""
$grid[18,4] = 7
$grid[18,5] = ""
$grid[18,6] = 1
$grid[18,7] = 6
$grid[18,8] = ""
$grid[18,9] = ""
$grid[18,10] = ""
$grid[18,11] = ""
$grid[18,12] = ""
$grid[18,13] = ""
$grid[18,14] = ""
$grid[18,15] = ""
$grid[18,16] = ""
$grid[18,17] = 2
$grid[18,18] = ""
$grid[18,19] = ""
$grid[18,20] = ""
$grid[19,0] = 8
$grid[19,1] = ""
$grid[19,2] = ""
$grid[19,3] = ""
$grid[19,4] = 1
$grid[19,5] = ""
$grid[19,6] = ""
$grid[19,7] = ""
$grid[19,8] = 7
$grid[19,9] = ""
$grid[19,10] = ""
$grid[19,11] = ""
$grid[19,12] = ""
$grid[19,13] = ""
$grid[19,14] = 2
$grid[19,15] = 9
$grid[19,16] = ""
$grid[19,17] = 4
$grid[19,18] = 5
$grid[19,19] = 8
$grid[19,20] = ""
$grid[20,0] = ""
$grid[20,1] = ""
$grid[20,2] = 9
$grid[20,3] = ""
$grid[20,4] = ""
$grid[20,5] = ""
$grid[20,6] = ""
$grid[20,7] = ""
$grid[20,8] = ""
$grid[20,9] = ""
$grid[20,10] = ""
$grid[20,11] = ""
$grid[20,12] = ""
$grid[20,13] = 1
$grid[20,14] = 5
$grid[20,15] = ""
$grid[20,16] = ""
$grid[20,17] = ""
$grid[20,18] = 4
$grid[20,19] = ""
$grid[20,20] = ""

$ws.Range("A1:U21").Value = $grid

$ws.Activate()
